# Generate Report for Handback
# Update the timestamp cells on the Overview, zh-cn and de-de sheets to
# reflect the latest handback / handoff / xliff-generation times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-30 23:11:24"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-30 23:11:19"
$wsZhCn.Range("K2").Value = "2016-08-30 23:11:37"

# de-de sheet: "Correspond Handoff Datetime" (shared with Overview's value)
# and "Correspond Handback DateTime" for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-30 23:11:24"
$wsDeDe.Range("K2").Value = "2016-08-30 23:11:45"
